# Apply the "BluetoothSpeakers-RyanAntolin" documentation update.
# This updates several DataEntry answers (B9, B10, B11, B12, B13, B14, B15),
# adjusts related row heights / a font style, and refreshes the selection.
# The DataBase sheet recalculates automatically because its cells are
# formulas that reference DataEntry (e.g. =DataEntry!B9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

# --- Row 9: collaborators added to the existing answer ---
$ws.Range("B9").Value = "Bluetooth, Alexa API(Amazon Voice) Patrick Ng, Ruel John Cooutauco"

# --- Row 10: "50 word problem statement" now has an answer ---
$ws.Range("B10").Value = "Bluetooth Sensor project is to be able to connect to a speaker via bluetooth connection wirelessly. How this works is that you take a mobile device and connect to the bluetooth device in order to produce sound waves from a speaker. This technology has been around for awhile but would want to recreate it."
$ws.Rows.Item(10).RowHeight = 75

# --- Row 11: "100 words of background" now has an answer ---
$ws.Range("B11").Value = "`nThis project is to connect from a mobile or Bluetooth connection that can play audio towards the Bluetooth speaker. This is so that it can play audio in a speaker form via Bluetooth connection. It will amplify the audio so it can play in areas that you usually cannot hear with your mobile speakers. It is more efficient than playing it through an auxiliary cord because it can play wirelessly via Bluetooth which would have less hassle with the auxiliary cord. Also the speaker can possibly have batteries on it so that it can have it be portable.`n"
$ws.Rows.Item(11).RowHeight = 150

# --- Row 15: "Solution description" now has an answer ---
$ws.Range("B15").Value = "To make bluetooth usable with audio and amplify the sound more in certain areas."
$ws.Rows.Item(15).RowHeight = 30

# --- Row 13: "Existing research IEEE paper APA citation" replaced ---
$ws.Range("B13").Value = "Lumpkins, W. (n.d.). The MobiAria Wireless Bluetooth Speaker. Retrieved September 18, 2017, from http://ieeexplore.ieee.org/document/6685931/"

# --- Row 12: "Current product APA citation" replaced ---
$ws.Range("B12").Value = "Bluetooth in wireless communication. (n.d.). Retrieved September 18, 2017, from http://ieeexplore.ieee.org/document/1007414/"

# --- Row 14: "Brief description of planned purchases" replaced, with a
#     distinct font (Times New Roman 12, dark grey) applied to the answer ---
$ws.Range("B14").WrapText = $false
$cell14Font = $ws.Range("B14").Font
$cell14Font.Name = "Times New Roman"
$cell14Font.Size = 12
$cell14Font.Color = 3355443
$cell14Font.Family = 1
$ws.Range("B14").Value = "Bluetooth in wireless communication. (n.d.). Retrieved September 18, 2017, from http://ieeexplore.ieee.org/document/1007414/"

# --- Update the saved selection / scroll position on the DataEntry sheet ---
$ws.Activate()
$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# Recalculate so DataBase's cached formula results reflect the new values.
$excel.CalculateFull()
